# Auto-generated: updates cryptos list Price (D) and Volume(1h) (E) columns
# to match the Sat Oct 14 13:25:53 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text cell (many of these price/volume
    # strings look numeric, e.g. '206.29' or '7.28', and Excel's
    # normal paste/typing heuristics would silently convert them to
    # numbers). Apply a Text number format before writing the value,
    # matching how a user would paste these as text, then restore the
    # cell to the workbook's default 'Normal' style so no visible
    # formatting changes are introduced.
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

# --- Price (column D) updates ---
Set-TextValue $ws.Range('D2') '26.930.79'
Set-TextValue $ws.Range('D3') '1.551.85'
Set-TextValue $ws.Range('D5') '206.29'
Set-TextValue $ws.Range('D10') '0.0594'
Set-TextValue $ws.Range('D12') '1.772.67'
Set-TextValue $ws.Range('D13') '1.552.27'
Set-TextValue $ws.Range('D16') '26.908.36'
Set-TextValue $ws.Range('D17') '61.62'
Set-TextValue $ws.Range('D18') '0.0₃0712'
Set-TextValue $ws.Range('D19') '217.03'
Set-TextValue $ws.Range('D20') '7.28'
Set-TextValue $ws.Range('D23') '9.19'
Set-TextValue $ws.Range('D25') '153.65'
Set-TextValue $ws.Range('D30') '0.0469'
Set-TextValue $ws.Range('D34') '1.409.98'
Set-TextValue $ws.Range('D36') '0.965'
Set-TextValue $ws.Range('D37') '2.29'
Set-TextValue $ws.Range('D39') '0.524'
Set-TextValue $ws.Range('D40') '0.807'
Set-TextValue $ws.Range('D45') '64.44'
Set-TextValue $ws.Range('D47') '1.686.56'
Set-TextValue $ws.Range('D48') '87.33'
Set-TextValue $ws.Range('D49') '0.0519'
Set-TextValue $ws.Range('D51') '0.0958'

# --- Volume(1h) (column E) updates ---
Set-TextValue $ws.Range('E2') '  -0.33%  '
Set-TextValue $ws.Range('E3') '  -0.42%  '
Set-TextValue $ws.Range('E4') '  -0.47%  '
Set-TextValue $ws.Range('E5') '  -0.57%  '
Set-TextValue $ws.Range('E6') '  +0.45%  '
Set-TextValue $ws.Range('E7') '  -0.46%  '
Set-TextValue $ws.Range('E8') '  +1.08%  '
Set-TextValue $ws.Range('E9') '  -0.48%  '
Set-TextValue $ws.Range('E10') '  +0.57%  '
Set-TextValue $ws.Range('E11') '  -0.76%  '
Set-TextValue $ws.Range('E12') '  -0.44%  '
Set-TextValue $ws.Range('E13') '  -0.38%  '
Set-TextValue $ws.Range('E14') '  +0.38%  '
Set-TextValue $ws.Range('E15') '  +0.21%  '
Set-TextValue $ws.Range('E17') '  -0.63%  '
Set-TextValue $ws.Range('E18') '  +3.34%  '
Set-TextValue $ws.Range('E19') '  +0.29%  '
Set-TextValue $ws.Range('E20') '  +0.00%  '
Set-TextValue $ws.Range('E22') '  +1.08%  '
Set-TextValue $ws.Range('E23') '  -0.46%  '
Set-TextValue $ws.Range('E25') '  +0.32%  '
Set-TextValue $ws.Range('E26') '  -0.38%  '
Set-TextValue $ws.Range('E27') '  +0.03%  '
Set-TextValue $ws.Range('E28') '  +0.33%  '
Set-TextValue $ws.Range('E29') '  -0.31%  '
Set-TextValue $ws.Range('E30') '  +1.06%  '
Set-TextValue $ws.Range('E31') '  -0.98%  '
Set-TextValue $ws.Range('E32') '  -0.60%  '
Set-TextValue $ws.Range('E33') '  +3.55%  '
Set-TextValue $ws.Range('E34') '  +0.65%  '
Set-TextValue $ws.Range('E35') '  +2.06%  '
Set-TextValue $ws.Range('E36') '  +0.39%  '
Set-TextValue $ws.Range('E37') '  +0.09%  '
Set-TextValue $ws.Range('E38') '  +0.18%  '
Set-TextValue $ws.Range('E39') '  +0.09%  '
Set-TextValue $ws.Range('E40') '  -0.57%  '
Set-TextValue $ws.Range('E41') '  -0.45%  '
Set-TextValue $ws.Range('E42') '  +2.98%  '
Set-TextValue $ws.Range('E43') '  +1.61%  '
Set-TextValue $ws.Range('E44') '  +0.46%  '
Set-TextValue $ws.Range('E45') '  +0.66%  '
Set-TextValue $ws.Range('E46') '  -1.28%  '
Set-TextValue $ws.Range('E47') '  -0.46%  '
Set-TextValue $ws.Range('E48') '  +1.25%  '
Set-TextValue $ws.Range('E49') '  +1.34%  '
Set-TextValue $ws.Range('E50') '  +2.95%  '
Set-TextValue $ws.Range('E51') '  -0.44%  '
